$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45432

$ws.Range("D29").Value = 547
$ws.Range("D30").Value = 547
$ws.Range("D31").Value = 547
$ws.Range("D32").Value = 547
$ws.Range("D33").Value = 1615
